$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the locator value in C5 (was the 4th leaflet-zoom-animated div xpath,
# now the "Zoom out" anchor xpath used by the new API methods/utils)
$ws.Range("C5").Value = "//a[@title='Zoom out']"

# Reflect the new active selection on the sheet
$ws.Range("C5").Select()
